# Update the "pipelines" sheet:
#  - Header C1: "jsonFile" -> "AppFolderPath"
#  - Row 2 (Pipeline 1) col C: "Apps\\app1.json" -> "Apps\\App1"
#  - Row 3 (Pipeline 2) col C: "Apps\\app2.json" -> "Apps\\App2"
# (literal double backslash, as stored in the original file).
# The "Apps" / "\\" prefix keeps its original rich-text coloring (two
# separate colored runs matching the VS-Code-like theme used elsewhere in
# the file); only the filename portion's text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("C1").Value = "AppFolderPath"

# --- Row 2: Pipeline 1 / App1 ---
# Note: the leading "Apps" run intentionally keeps the default (unset) font,
# matching the original file's formatting - only the "\\" (two literal
# backslashes, unchanged by this edit) and filename runs carry explicit
# Consolas/coloured rPr.
$ws.Range("C2").Value = "Apps\\App1"

$r2b = $ws.Range("C2").Characters(5, 2)      # "\\"
$r2b.Font.Name = "Consolas"
$r2b.Font.Size = 11
$r2b.Font.Color = 8239831                    # FFD7BA7D (R=0xD7,G=0xBA,B=0x7D)

$r2c = $ws.Range("C2").Characters(7, 4)      # "App1"
$r2c.Font.Name = "Consolas"
$r2c.Font.Size = 11
$r2c.Font.Color = 7901646                    # FFCE9178 (R=0xCE,G=0x91,B=0x78)

# --- Row 3: Pipeline 2 / App2 ---
$ws.Range("C3").Value = "Apps\\App2"

$r3b = $ws.Range("C3").Characters(5, 2)      # "\\"
$r3b.Font.Name = "Consolas"
$r3b.Font.Size = 11
$r3b.Font.Color = 8239831                    # FFD7BA7D

$r3c = $ws.Range("C3").Characters(7, 4)      # "App2"
$r3c.Font.Name = "Consolas"
$r3c.Font.Size = 11
$r3c.Font.Color = 7901646                    # FFCE9178

# --- Selection state: C3 selected first, F4 made the active cell ---
$excel.Union($ws.Range("C3"), $ws.Range("F4")).Select()

Write-Output "edit applied"
